$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("ZoneLetter") holds values "T" for data rows 2-180; change them to "V".
for ($r = 2; $r -le 180; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "T") {
        $cell.Value = "V"
    }
}
